# "Generate Report for Handoff"
# Updates the localization-status report: the Overview/zh-cn/de-de sheets move
# from "Handed back: in sync with en-US" to "Ready for handoff", and the
# corresponding handoff timestamps are refreshed. The two narrow status/date
# columns on each sheet are then re-sized to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-25 15:03:29"

# ---- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-25 15:03:25"

# ---- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-25 15:03:29"

# ---- Shrink the status columns to fit the new, shorter text ------------
# (Excel auto-fits the "zh-cn"/"de-de" status columns on the Overview sheet
# and the "Status" column on each language sheet after the text changes.)
$overview.Columns.Item(5).ColumnWidth = 16.3   # column E ("zh-cn")
$overview.Columns.Item(6).ColumnWidth = 16.3   # column F ("de-de")
$zhcn.Columns.Item(3).ColumnWidth = 16.3        # column C ("Status")
$dede.Columns.Item(3).ColumnWidth = 16.3        # column C ("Status")
